$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 113.1
$ws.Range("I5").Value = 126
$ws.Range("K5").Value = 126
$ws.Range("M5").Value = -11
$ws.Range("H18").Value = 333
$ws.Range("I18").Value = 399.5
$ws.Range("J18").Value = 200
$ws.Range("K18").Value = 399.5
$ws.Range("L18").Value = 200
$ws.Range("M18").Value = -115.5
$ws.Range("N18").Value = -768
$ws.Range("H51").Value = 6688.6665
$ws.Range("I51").Value = 5037.4287
$ws.Range("J51").Value = 9000.4
$ws.Range("K51").Value = 5037.4287
$ws.Range("L51").Value = 9000.4
$ws.Range("M51").Value = -4553.4287
$ws.Range("N51").Value = -9968.4
$ws.Range("H127").Value = 11510.223
$ws.Range("I127").Value = 1644.0667
$ws.Range("K127").Value = 4932.2001
$ws.Range("M127").Value = 27.79989999999998
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9320.104499999999
$ws.Range("I32").Value = 7028.4424
$ws.Range("K32").Value = 7028.4424
$ws.Range("M32").Value = -6741.4424
$ws.Range("H45").Value = 5150
$ws.Range("I45").Value = 3436.5715
$ws.Range("K45").Value = 3436.5715
$ws.Range("M45").Value = -3059.5715
$ws.Range("H97").Value = 7938728.5
$ws.Range("I97").Value = 2273.875
$ws.Range("K97").Value = 2273.875
$ws.Range("M97").Value = -1777.875
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H122").Value = 3974.825
$ws.Range("I122").Value = 2556.2
$ws.Range("J122").Value = 5393.45
$ws.Range("K122").Value = 7668.599999999999
$ws.Range("L122").Value = 16180.35
$ws.Range("M122").Value = -5218.599999999999
$ws.Range("N122").Value = -21080.35
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3632.3
$ws.Range("I20").Value = 3036
$ws.Range("J20").Value = 8999
$ws.Range("K20").Value = 3036
$ws.Range("L20").Value = 8999
$ws.Range("M20").Value = -2789
$ws.Range("N20").Value = -9493
$ws.Range("H107").Value = 793.3333
$ws.Range("I107").Value = 793.3333
$ws.Range("K107").Value = 793.3333
$ws.Range("M107").Value = 1126.6667
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 335
$ws.Range("I5").Value = 412.75
$ws.Range("K5").Value = 412.75
$ws.Range("M5").Value = -300.75
$ws.Range("H7").Value = 280.47058
$ws.Range("J7").Value = 332.7143
$ws.Range("L7").Value = 332.7143
$ws.Range("N7").Value = -558.7143
$ws.Range("H22").Value = 5000
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H86").Value = 29395.166
$ws.Range("I86").Value = 41670
$ws.Range("J86").Value = 4845.5
$ws.Range("K86").Value = 41670
$ws.Range("L86").Value = 4845.5
$ws.Range("M86").Value = -40547
$ws.Range("N86").Value = -7091.5
$ws.Range("H89").Value = 29395.166
$ws.Range("I89").Value = 41670
$ws.Range("J89").Value = 4845.5
$ws.Range("K89").Value = 208350
$ws.Range("L89").Value = 24227.5
$ws.Range("M89").Value = -202734
$ws.Range("N89").Value = -35459.5
$ws.Range("H107").Value = 6050.2104
$ws.Range("J107").Value = 8238.77
$ws.Range("L107").Value = 8238.77
$ws.Range("N107").Value = -12078.77
$ws.Range("H134").Value = 3729.24
$ws.Range("I134").Value = 2560.3235
$ws.Range("K134").Value = 7680.970499999999
$ws.Range("M134").Value = -5145.970499999999
$ws.Range("H141").Value = 92162.39999999999
$ws.Range("J141").Value = 92162.39999999999
$ws.Range("L141").Value = 92162.39999999999
$ws.Range("N141").Value = -102522.4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 4034.5
$ws.Range("I51").Value = 4034.5
$ws.Range("K51").Value = 12103.5
$ws.Range("M51").Value = -11643.5
$ws.Range("H61").Value = 50
$ws.Range("I61").Value = 50
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 150
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = 65
$ws.Range("N61").ClearContents()
$ws.Range("H121").Value = 1933.1515
$ws.Range("J121").Value = 2306.6667
$ws.Range("L121").Value = 6920.000100000001
$ws.Range("N121").Value = -9540.000100000001
$ws.Range("H131").Value = 5747.5
$ws.Range("I131").Value = 3711.4285
$ws.Range("J131").Value = 20000
$ws.Range("K131").Value = 11134.2855
$ws.Range("L131").Value = 60000
$ws.Range("M131").Value = -6094.2855
$ws.Range("N131").Value = -70080
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 9994.200000000001
$ws.Range("J47").Value = 9994.5
$ws.Range("L47").Value = 9994.5
$ws.Range("N47").Value = -11130.5
$ws.Range("H49").Value = 9795
$ws.Range("I49").Value = 9795
$ws.Range("K49").Value = 9795
$ws.Range("M49").Value = -9611
$ws.Range("H70").Value = 228200
$ws.Range("I70").Value = 280250
$ws.Range("J70").Value = 20000
$ws.Range("K70").Value = 280250
$ws.Range("L70").Value = 20000
$ws.Range("M70").Value = -279980
$ws.Range("N70").Value = -20540
$ws.Range("H73").Value = 228200
$ws.Range("I73").Value = 280250
$ws.Range("J73").Value = 20000
$ws.Range("K73").Value = 280250
$ws.Range("L73").Value = 20000
$ws.Range("M73").Value = -279314
$ws.Range("N73").Value = -21872
$ws.Range("H86").Value = 120000
$ws.Range("J86").Value = 120000
$ws.Range("L86").Value = 120000
$ws.Range("N86").Value = -122372
$ws.Range("H89").Value = 120000
$ws.Range("J89").Value = 120000
$ws.Range("L89").Value = 360000
$ws.Range("N89").Value = -371856
$ws.Range("H107").Value = 916.6667
$ws.Range("I107").Value = 900
$ws.Range("J107").Value = 950
$ws.Range("K107").Value = 900
$ws.Range("L107").Value = 950
$ws.Range("M107").Value = 1020
$ws.Range("N107").Value = -4790
$ws.Range("H113").Value = 7249.9287
$ws.Range("I113").Value = 4288.778
$ws.Range("K113").Value = 4288.778
$ws.Range("M113").Value = -2118.778
$ws.Range("H119").Value = 42920
$ws.Range("J119").Value = 42920
$ws.Range("L119").Value = 42920
$ws.Range("N119").Value = -52596
$ws.Range("H122").Value = 3025.9062
$ws.Range("I122").Value = 1500.9412
$ws.Range("K122").Value = 4502.8236
$ws.Range("M122").Value = -2052.8236
$ws.Range("H132").Value = 3244.4102
$ws.Range("I132").Value = 2347.2
$ws.Range("K132").Value = 7041.599999999999
$ws.Range("M132").Value = -4511.599999999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("H96").Value = 2219
$ws.Range("I96").Value = 2154.6667
$ws.Range("J96").Value = 2283.3333
$ws.Range("K96").Value = 2154.6667
$ws.Range("L96").Value = 2283.3333
$ws.Range("M96").Value = -781.6667000000002
$ws.Range("N96").Value = -5029.3333
$ws.Range("H122").Value = 3771.4583
$ws.Range("I122").Value = 2588.4375
$ws.Range("J122").Value = 6137.5
$ws.Range("K122").Value = 7765.3125
$ws.Range("L122").Value = 18412.5
$ws.Range("M122").Value = -5315.3125
$ws.Range("N122").Value = -23312.5
$ws.Range("H132").Value = 2572.6843
$ws.Range("I132").Value = 1977.3636
$ws.Range("J132").Value = 3391.25
$ws.Range("K132").Value = 5932.0908
$ws.Range("L132").Value = 10173.75
$ws.Range("M132").Value = -3402.0908
$ws.Range("N132").Value = -15233.75
